$wb = $excel.ActiveWorkbook

# --- Relations sheet: swap source/target for relation "s" (row 4: C4/D4) ---
$relations = $wb.Worksheets.Item("Relations")
$c4 = $relations.Range("C4").Value2
$d4 = $relations.Range("D4").Value2
$relations.Range("C4").Value = $d4
$relations.Range("D4").Value = $c4

# --- Terms sheet: update the displayed term text for t7 ("s[A*B]" -> "s[B*A]") ---
$terms = $wb.Worksheets.Item("Terms")
$terms.Range("B10").Value = "s[B*A]"

# --- Update the active sheet / selections to match the saved view state ---
# Active sheet moves from Relations to Terms, with a new selection on Terms.
$terms.Activate()
$terms.Range("C13").Select()
